$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (they contain dotted/grouped numeric-looking strings
# like "66.462.96" or plain decimals like "0.127" that Excel would otherwise coerce to numbers).
$dCells = @("D2","D3","D5","D6","D9","D12","D14","D15","D16","D17","D18","D19","D20","D21","D22","D25","D29","D30","D32","D34","D35","D36","D37","D40","D41","D43","D44","D45","D47","D48","D49")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.462.96"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.320.98"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "586.63"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("D6").Value = "183.01"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +8.11%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.127"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "3.897.47"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  -4.37%  "
$ws.Range("D14").Value = "66.461.46"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "26.45"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.318.20"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "431.39"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "13.34"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("D20").Value = "5.55"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "7.46"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("D22").Value = "72.33"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "3.444.97"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("E28").Value = "  -3.28%  "
$ws.Range("D29").Value = "9.03"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "22.49"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "5.24"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "1.21"
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "6.64"
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("D37").Value = "159.62"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.887.81"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "26.86"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").Value = "4.35"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "40.35"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "0.0669"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").Value = "2.33"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("D48").Value = "23.45"
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("D49").Value = "318.38"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  +4.29%  "
